$d = $word.ActiveDocument

# --- Step 1: fix the run split in paragraph 2 around the "_GoBack" bookmark ---
# Originally: "...Operation Mapping " [bookmark] "parameter" " list belong to operation mapping $" "Operation" "_Name" "."
# Target:     "...Operation Mapping " "parameter" " list belong to operation mapping " " " [bookmark] "$" "Operation" "_Name" "."

# 1a. Collapse the two runs after the bookmark ("parameter" + " list belong to operation mapping $")
#     down to a single "$" run.
$d.Content.Find.Execute("parameter list belong to operation mapping `$", $true, $false, $false, $false, $false, $true, 1, $false, "`$", 1) | Out-Null

# 1b. Insert the now-missing text right before the bookmark (so it reads identically to before,
#     but now lives ahead of the bookmark instead of after it), plus the extra trailing space.
$bm = $d.Bookmarks("_GoBack")
$insBm = $d.Range($bm.Start, $bm.Start)
$insBm.InsertBefore("parameter list belong to operation mapping  ")

# --- Step 2: append the new "Operation Mapping:" heading + details table ---
$p2 = $d.Paragraphs(2)
$insEnd = $d.Range($p2.Range.End, $p2.Range.End)

$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Operation Mapping: </w:t>
  </w:r>
</w:p>
<w:tbl>
  <w:tblPr>
    <w:tblW w:w="5000" w:type="pct"/>
    <w:tblBorders>
      <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
      <w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>
      <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
      <w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>
      <w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/>
      <w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/>
    </w:tblBorders>
    <w:tblLook w:val="01E0" w:firstRow="1" w:lastRow="1" w:firstColumn="1" w:lastColumn="1" w:noHBand="0" w:noVBand="0"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="3145"/>
    <w:gridCol w:w="6205"/>
  </w:tblGrid>
  <w:tr>
    <w:trPr>
      <w:trHeight w:val="314"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="5000" w:type="pct"/>
        <w:gridSpan w:val="2"/>
        <w:shd w:val="clear" w:color="auto" w:fill="8496B0" w:themeFill="text2" w:themeFillTint="99"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="22"/>
            <w:szCs w:val="22"/>
            <w:lang w:eastAsia="zh-CN"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="22"/>
            <w:szCs w:val="22"/>
            <w:lang w:eastAsia="zh-CN"/>
          </w:rPr>
          <w:t>Operation Mapping</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:trHeight w:val="314"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1682" w:type="pct"/>
        <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="22"/>
            <w:szCs w:val="22"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="22"/>
            <w:szCs w:val="22"/>
          </w:rPr>
          <w:t>Name</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="3318" w:type="pct"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t>$Name_Value</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1682" w:type="pct"/>
        <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="22"/>
            <w:szCs w:val="22"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="22"/>
            <w:szCs w:val="22"/>
          </w:rPr>
          <w:t>NameSpace</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="3318" w:type="pct"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t>$NameSpace_Value</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1682" w:type="pct"/>
        <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="22"/>
            <w:szCs w:val="22"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
            <w:b/>
            <w:sz w:val="22"/>
            <w:szCs w:val="22"/>
          </w:rPr>
          <w:t>Software component version</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="3318" w:type="pct"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t>$ComponentID_Value</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insEnd.InsertXML($xmlSnippet)
